$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.491.44"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.672.33"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.671.06"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.94%  "
$ws.Range("E10").Value = "  +5.00%  "
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.05"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.161.76"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "72.358.69"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.43"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.674.56"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.94"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "378.75"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.20"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.53%  "
$ws.Range("E23").Value = "  +10.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.46"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.37"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.94"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.807.10"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.14"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "521.39"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.62"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.08%  "
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("E40").Value = "  -7.11%  "
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.06"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.335"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.31"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.43"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.549"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("E51").Value = "  +2.12%  "
